$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 28
$ws.Range("H4").Value = 22
$ws.Range("H6").Value = 36
$ws.Range("H8").Value = 4
$ws.Range("H9").Value = 52

$ws.Range("H9").Select()
